# Insert a new "count_ratio_dict" worksheet right after "quota_dict" (i.e. before
# "proposal_dictionary"), matching the new sheet order:
#   ... quota_dict, count_ratio_dict, proposal_dictionary, eg_colors,
#       basic_job_colors, enhanced_job_colors
$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("quota_dict")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "count_ratio_dict"

# Header row
$headers = @("basic_job", "group1", "group2", "group3", "weight1", "weight2", "weight3", "cap", "month_start", "month_end")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# Data row 2
$row2 = @(1, 1, "2, 3 ", 0, 2.48, 1, 0, 320, 34, 65)
for ($col = 1; $col -le $row2.Length; $col++) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Value = $row2[$col - 1]
    $cell.HorizontalAlignment = -4108
}

# Data row 3
$row3 = @(4, 1, 2, 3, 2.46, 1, 1.2, 580, 34, 55)
for ($col = 1; $col -le $row3.Length; $col++) {
    $cell = $ws.Cells.Item(3, $col)
    $cell.Value = $row3[$col - 1]
    $cell.HorizontalAlignment = -4108
}

# Make this the active sheet / selection, matching the new workbook's activeTab
[void]$ws.Activate()
[void]$ws.Range("A4").Select()
